# Weekly update: a new "Perejil" (Vega Monumental Concepción) price report pair
# (Primera / Segunda quality rows) is published, pushing the whole history
# table down by two rows. We reproduce this by duplicating the current first
# data pair (rows 176:177) one slot down via Copy+Insert (which shifts every
# following row down by two, spilling the former last pair into two brand new
# rows 230:231), then stamping the newly freed top rows with the new report
# date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate rows 176:177 and insert the copy in place - this shifts rows
# 176:229 down to 178:231, preserving all values/styles, and grows the used
# range from R229 to R231.
$ws.Rows("176:177").Copy()
$ws.Rows("176:177").Insert()

# The newly inserted top rows represent the latest weekly report date.
$ws.Range("D176").Value2 = 45120
$ws.Range("D177").Value2 = 45120
